# The presentation ships with two theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" colours (currently linked from the Notes Master)
#   ppt/theme/theme2.xml -> "Integral" colours      (currently linked from the Slide Master,
#                                                     i.e. the design actually applied to the deck)
#
# The authored change swaps the content of the two theme parts: the design that drives the
# slide master/layouts/slides becomes the default "Office Theme" colour palette, while the
# "Integral" palette moves over to what used to hold the Office Theme.
#
# This host only exposes a single live theme object (reached from SlideMaster / NotesMaster /
# HandoutMaster alike) and persists edits to it back into the theme part referenced by the
# Slide Master. We reproduce the swap by rewriting each of the twelve theme colour slots
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) from the current "Integral" values over to
# the "Office Theme" values, using the documented mutation path: ThemeColorScheme.Colors(i).RGB.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$scheme = $master.Theme.ThemeColorScheme

$scheme.Colors(1).RGB = 0          # dk1      -> 000000
$scheme.Colors(2).RGB = 16777215   # lt1      -> FFFFFF
$scheme.Colors(3).RGB = 6968388    # dk2      -> 44546A
$scheme.Colors(4).RGB = 15132391   # lt2      -> E7E6E6
$scheme.Colors(5).RGB = 13998939   # accent1  -> 5B9BD5
$scheme.Colors(6).RGB = 3243501    # accent2  -> ED7D31
$scheme.Colors(7).RGB = 10855845   # accent3  -> A5A5A5
$scheme.Colors(8).RGB = 49407      # accent4  -> FFC000
$scheme.Colors(9).RGB = 12874308   # accent5  -> 4472C4
$scheme.Colors(10).RGB = 4697456   # accent6  -> 70AD47
$scheme.Colors(11).RGB = 12673797  # hlink    -> 0563C1
$scheme.Colors(12).RGB = 7491477   # folHlink -> 954F72
